# loading in full final dataset, preprocessing/data cleaning, save climate
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data cleaning: species in A15 was misidentified; update to correct species.
$ws.Range("A15").Value = "Lupinus bicolor"

# Move the active selection to the corrected cell.
$ws.Range("A15").Select()
